$wb = $excel.ActiveWorkbook

# --- Sheet1: update/add asset price history values ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A1").Value = "6910,7"
$ws1.Range("A2").Value = "1495,21"
$ws1.Range("A4").Value = "12.23"
$ws1.Range("A10").Value = "2,199542"
$ws1.Range("A34").Value = "60,81"

# --- data sheet: new asset (ethereum), sheet/cell reference, and currency (PLN) ---
$wsd = $wb.Worksheets.Item("data")
$wsd.Range("A1").Value = "ethereum"
$wsd.Range("B1").Value = ""
$wsd.Range("C1").Value = ""
$wsd.Range("A2").Value = "Sheet1"
$wsd.Range("B2").Value = ""
$wsd.Range("A3").Value = "A1"
$wsd.Range("B3").Value = ""
$wsd.Range("A4").Value = "PLN"
